# Fixed spelling mistakes ("treament" -> "treatment") and added "_id"/".id"
# suffixes to the machine-readable field names stored in the hidden row 3
# of the sideload template (row 1 = section headers, row 2 = human-readable
# headers, row 3 = machine field names used for CSV/JSON export).
#
# Row 2 (human readable headers) is left untouched; only the hidden row 3
# values change:
#   organism            -> organism_id
#   experiment_type     -> experiment_type_id
#   instrument          -> instrument_id
#   proteomic_fraction  -> proteomic_fraction_id
#   sample_type         -> sample_type_id
#   cell_type           -> cell_type_id
#   treament.L.inhibitor.name          -> treatment.L.inhibitor.id
#   treament.L.inhibitor.concentration -> treatment.L.inhibitor.concentration
#   treament.L.inhibitor.time          -> treatment.L.inhibitor.time
#   treament.L.probe.name              -> treatment.L.probe.id
#   treament.L.probe.concentration     -> treatment.L.probe.concentration
#   treament.L.probe.time              -> treatment.L.probe.time
#   treament.H.inhibitor.name          -> treatment.H.inhibitor.id
#   treament.H.inhibitor.concentration -> treatment.H.inhibitor.concentration
#   treament.H.inhibitor.time          -> treatment.H.inhibitor.time
#   treament.H.probe.name              -> treatment.H.probe.id
#   treament.H.probe.concentration     -> treatment.H.probe.concentration
#   treament.H.probe.time              -> treatment.H.probe.time
#
# (the ".method" fields already used the correct "treatment." spelling and
# are left untouched)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D3").Value2 = "organism_id"
$ws.Range("E3").Value2 = "experiment_type_id"
$ws.Range("F3").Value2 = "instrument_id"
$ws.Range("G3").Value2 = "proteomic_fraction_id"
$ws.Range("H3").Value2 = "sample_type_id"
$ws.Range("I3").Value2 = "cell_type_id"

$ws.Range("N3").Value2 = "treatment.L.inhibitor.time"
$ws.Range("P3").Value2 = "treatment.L.probe.concentration"
$ws.Range("O3").Value2 = "treatment.L.probe.id"
$ws.Range("R3").Value2 = "treatment.L.probe.time"
$ws.Range("S3").Value2 = "treatment.H.inhibitor.id"
$ws.Range("K3").Value2 = "treatment.L.inhibitor.id"
$ws.Range("L3").Value2 = "treatment.L.inhibitor.concentration"
$ws.Range("T3").Value2 = "treatment.H.inhibitor.concentration"
$ws.Range("V3").Value2 = "treatment.H.inhibitor.time"
$ws.Range("X3").Value2 = "treatment.H.probe.concentration"
$ws.Range("Z3").Value2 = "treatment.H.probe.time"
$ws.Range("W3").Value2 = "treatment.H.probe.id"

# Re-fit the (hidden) row after editing its contents so no spurious explicit
# row-height override gets persisted.
$ws.Rows.Item(3).AutoFit()
